# Insert a new record row at row 17, shifting the existing rows 17-45 down to 18-46.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Insert()

$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value = "Ñuble"
$ws.Cells.Item(17, 4).Value = 44690
$ws.Cells.Item(17, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 6).Value = 100112040
$ws.Cells.Item(17, 7).Value = "Cilantro"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 120
$ws.Cells.Item(17, 11).Value = 550
$ws.Cells.Item(17, 12).Value = 600
$ws.Cells.Item(17, 13).Value = 575
$ws.Cells.Item(17, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(17, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(17, 16).Value = 575
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = "Hortaliza"
